# Creating the railway line of the overland
# 1. Rename a handful of station names on "Stations" to add the " station" suffix.
# 2. Insert a new "Overland Line" worksheet between "Stations" and "TrainServices"
#    listing the stations (and their order) that the Overland service travels through.

$wb = $excel.ActiveWorkbook

$stations = $wb.Worksheets.Item("Stations")
$trainServices = $wb.Worksheets.Item("TrainServices")

# --- 1. Rename a few station names (append "station") ---------------------
$stations.Range("B7").Value  = "Yass Junction station"
$stations.Range("B30").Value = "Broken Hill station"
$stations.Range("B31").Value = "Hornsby station"
$stations.Range("B32").Value = "Gosford station"
$stations.Range("B33").Value = "Wyong station"
$stations.Range("B34").Value = "Fassifern station"

# Reflect the scrolled/selected state the author ended up with on the Stations tab.
$stations.Activate()
$stations.Range("B86:B94").Select()

# --- 2. Insert the new "Overland Line" worksheet ---------------------------
$overland = $wb.Worksheets.Add($null, $stations)
$overland.Name = "Overland Line"

$overland.Columns.Item(1).ColumnWidth = 22.1666666666667

$overland.Range("A1").Value = "Station_name"
$overland.Range("B1").Value = "Order"

$stationOrder = @(
    @("Melbourne Southern Cross", 1),
    @("North Shore", 2),
    @("Arrarat", 3),
    @("Stawell", 4),
    @("Horsham", 5),
    @("Dimboola", 6),
    @("Nhill", 7),
    @("Adelaide Parklands", 10),
    @("Murray Bridge", 9),
    @("Bordertown", 8)
)

$row = 2
foreach ($item in $stationOrder) {
    $overland.Cells.Item($row, 1).Value = $item[0]
    $overland.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Make "Overland Line" the active tab/selection, matching the saved state.
$overland.Activate()
$overland.Range("D17").Select()
